$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8165858387947083
$ws.Range("B1").Value = 1.907322525978088
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.948338508605957
$ws.Range("E1").Value = 0.5562154054641724
